$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 now carries what used to be row 4's data
$ws.Range("A2").Value = "Você"
$ws.Range("B2").Value = "Desconhecido"
$ws.Range("C2").Value = "R$ 200,00"
$ws.Range("D2").Value = "COMERCIO DE POLPAS SOUZA E DIAS LTD..."
$ws.Range("E2").Value = "Funcionário"

# Remove the now-obsolete rows 3-5
$ws.Range("A3:E5").Delete()

# Narrow columns A and B
$ws.Columns.Item(1).ColumnWidth = 5.17
$ws.Columns.Item(2).ColumnWidth = 13.17
